# Matningsdata.xlsx — remove the "Fakturtotal (kr):" totals row.
#
# The source workbook had a trailing totals row (row 7: label in A7,
# the (currently zero) invoice total in B7) that the Kostnadsandel
# column (F2:F6) multiplied against via the absolute reference B$7.
# That row was removed from the sheet (and therefore from the
# Tabell3 table range, which shrinks from A1:F7 to A1:F6); the
# formulas that used to read B$7 now carry a dangling #REF! because
# their anchor row is gone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Deleting the entire row shifts everything below it up, updates the
# worksheet dimension, auto-shrinks the Tabell3 ListObject/AutoFilter
# range from A1:F7 to A1:F6, and rewrites every formula that pointed
# at a cell in row 7 (B$7 in F2:F6) into a #REF! error — exactly like
# using Home > Delete > Delete Sheet Rows on the row 7 selection.
$ws.Rows(7).Delete() | Out-Null

# Leave the selection where the author's last save left it.
$ws.Range("G4").Select() | Out-Null
